$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts old rows 3..25 down to 4..26),
# then copy the formatting (thin box border) from row 2 onto the new row.
$ws.Rows("3").Insert()
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)

# Populate the newly inserted row 3 with the new product line.
$ws.Range("A3").Value = "20141063"
$ws.Range("B3").Value = "SQ CSHW BTRFLY 3X22G"
$ws.Range("C3").Value = "SWASH"
$ws.Range("D3").Value = "1"
$ws.Range("E3").Value = "2"
$ws.Range("F3").Value = "TG,(E-1B)"

# The row that used to be row 3 (S/Q CHOCO CASHEW 52G) is now row 4;
# its sequence number changes from 2 to 3.
$ws.Range("E4").Value = "3"

# Fix the CADBURY product names (57 -> 52) on what are now rows 12 and 13.
$ws.Range("B12").Value = "CADBURY DAIRY MLK 52"
$ws.Range("B13").Value = "CADBURY CSHW NUT 52G"

Write-Output "done"
